$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pp.xpt")

# --- Row 6 ---
$ws.Range("A6").Value = "Sample001"
$ws.Range("B6").Value = "PP"
$ws.Range("C6").Value = "01-005"
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = "Part A Period 1"
$ws.Range("F6").Value = "AUCIFP"
$ws.Range("G6").Value = "AURC Infinity Obs Norm by Dose"
$ws.Range("H6").Value = "PPCAT01"
$ws.Range("I6").Value = 1194.546
$ws.Range("J6").Value = "day*ug/mL/mg"
$ws.Range("K6").Value = 1195
$ws.Range("L6").Value = 1195
$ws.Range("M6").Value = "day*ug/mL/mg"
$ws.Range("N6").Value = "PLASMA"
$ws.Range("O6").Value = "2018-04-09T09:05"

# --- Row 7 ---
$ws.Range("A7").Value = "Sample001"
$ws.Range("B7").Value = "PP"
$ws.Range("C7").Value = "01-006"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = "Part A Period 1"
$ws.Range("F7").Value = "AUCIFP"
$ws.Range("G7").Value = "AURC Infinity Obs Norm by Dose"
$ws.Range("H7").Value = "PPCAT01"
$ws.Range("I7").Value = 1194.546
$ws.Range("J7").Value = "h*ng/mL"
$ws.Range("K7").Value = 1195
$ws.Range("L7").Value = 1195
$ws.Range("M7").Value = "h*ng/mL"
$ws.Range("N7").Value = "PLASMA"
$ws.Range("O7").Value = "2018-04-09T09:05"

# --- Styles: copy formatting from row 5 down to rows 6 and 7 ---
$ws.Range("A5:O5").Copy()
$ws.Range("A6:O7").PasteSpecial(-4122) # xlPasteFormats

# --- Number format for L6/L7 ---
$ws.Range("L6:L7").NumberFormat = ".00"

# --- Column width for column G (closest achievable to 14.33203125 given engine's width quantization) ---
$ws.Columns.Item(7).ColumnWidth = 13.43

# --- Selection ---
$ws.Range("K7").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Workbook view xWindow (window restored to x=0) ---
$wb.Windows.Item(1).Left = 0

